$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.411.08"
$ws.Range("E2").Value = "  -3.93%  "
$ws.Range("D3").Value = "2.460.53"
$ws.Range("E3").Value = "  -6.78%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'549.91"
$ws.Range("E5").Value = "  -5.17%  "
$ws.Range("D6").Value = "'146.81"
$ws.Range("E6").Value = "  -6.65%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.593"
$ws.Range("E8").Value = "  -6.29%  "
$ws.Range("D9").Value = "2.459.78"
$ws.Range("E9").Value = "  -6.73%  "
$ws.Range("E10").Value = "  -9.77%  "
$ws.Range("E11").Value = "  -6.64%  "
$ws.Range("E12").Value = "  -1.89%  "
$ws.Range("D13").Value = "'0.352"
$ws.Range("E13").Value = "  -8.55%  "
$ws.Range("D14").Value = "'26.07"
$ws.Range("E14").Value = "  -9.65%  "
$ws.Range("D15").Value = "2.901.89"
$ws.Range("E15").Value = "  -6.78%  "
$ws.Range("E16").Value = "  -9.46%  "
$ws.Range("D17").Value = "61.315.72"
$ws.Range("E17").Value = "  -3.94%  "
$ws.Range("D18").Value = "2.458.78"
$ws.Range("E18").Value = "  -6.45%  "
$ws.Range("D20").Value = "'7.06"
$ws.Range("E20").Value = "  -9.13%  "
$ws.Range("E21").Value = "  -7.63%  "
$ws.Range("D22").Value = "'318.60"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("E24").Value = "  +0.32%  "
$ws.Range("D25").Value = "'64.04"
$ws.Range("E25").Value = "  -6.34%  "
$ws.Range("D26").Value = "0.0₃0980"
$ws.Range("E26").Value = "  -13.31%  "
$ws.Range("D27").Value = "'551.89"
$ws.Range("E27").Value = "  -5.79%  "
$ws.Range("D28").Value = "2.583.38"
$ws.Range("E28").Value = "  -6.58%  "
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("D30").Value = "'1.48"
$ws.Range("E30").Value = "  -10.90%  "
$ws.Range("E31").Value = "  -10.69%  "
$ws.Range("D32").Value = "'7.68"
$ws.Range("E32").Value = "  -6.82%  "
$ws.Range("E33").Value = "  -9.54%  "
$ws.Range("D34").Value = "'1.89"
$ws.Range("E34").Value = "  -8.37%  "
$ws.Range("E35").Value = "  -8.87%  "
$ws.Range("D36").Value = "'5.89"
$ws.Range("E36").Value = "  -11.62%  "
$ws.Range("E37").Value = "  -11.99%  "
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("E39").Value = "  -6.15%  "
$ws.Range("D40").Value = "'18.43"
$ws.Range("E40").Value = "  -7.01%  "
$ws.Range("E41").Value = "  -8.04%  "
$ws.Range("D42").Value = "'142.56"
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").Value = "'40.47"
$ws.Range("E44").Value = "  -4.10%  "
$ws.Range("D45").Value = "'2.37"
$ws.Range("E45").Value = "  -8.49%  "
$ws.Range("D46").Value = "'146.22"
$ws.Range("E46").Value = "  -10.07%  "
$ws.Range("E47").Value = "  -8.22%  "
$ws.Range("D48").Value = "'21.53"
$ws.Range("E48").Value = "  -11.08%  "
$ws.Range("E49").Value = "  -9.01%  "
$ws.Range("D50").Value = "'0.589"
$ws.Range("E50").Value = "  -7.06%  "
$ws.Range("D51").Value = "'0.0939"
$ws.Range("E51").Value = "  -6.81%  "
